# Insert a new weekly price record as row 149 on the "Cereza" sheet,
# pushing the existing rows 149-215 down to 150-216.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 149..215 down by inserting a brand-new blank row at 149.
$ws.Rows("149").Insert()

# Populate the newly inserted row with the new "Lapins" / "Primera" record.
$ws.Range("A149").Value = 10
$ws.Range("B149").Value = "Vega Modelo de Temuco"
$ws.Range("C149").Value = "La Araucanía"
$ws.Range("D149").Value = 44582
$ws.Range("E149").Value = 9
$ws.Range("F149").Value = "Fruta"
$ws.Range("G149").Value = 100103
$ws.Range("H149").Value = "Frutos de hueso (carozo)"
$ws.Range("I149").Value = 100103001
$ws.Range("J149").Value = "Cereza"
$ws.Range("K149").Value = "Lapins"
$ws.Range("L149").Value = "Primera"
$ws.Range("M149").Value = 295
$ws.Range("N149").Value = 11000
$ws.Range("O149").Value = 12000
$ws.Range("P149").Value = 11424
$ws.Range("Q149").Value = "`$/caja 18 kilos"
$ws.Range("R149").Value = "Región del Maule"
$ws.Range("S149").Value = 635
$ws.Range("T149").Value = 18

# Match the date number format used by the rest of the "Fecha" column.
$ws.Range("D149").NumberFormat = $ws.Range("D150").NumberFormat
